$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are stored as text, matching original inlineStr formatting,
# since some values look numeric (e.g. "0.999") and Excel would otherwise convert them.

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "60.890.95"
$ws.Cells.Item(2, 5).Value = "  -2.32%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.421.39"
$ws.Cells.Item(3, 5).Value = "  -1.05%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.03%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "569.34"
$ws.Cells.Item(5, 5).Value = "  -2.40%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "139.42"
$ws.Cells.Item(6, 5).Value = "  -3.33%  "

$ws.Cells.Item(7, 5).Value = "  +0.35%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.525"
$ws.Cells.Item(8, 5).Value = "  -1.14%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "2.402.80"
$ws.Cells.Item(9, 5).Value = "  -1.73%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.107"
$ws.Cells.Item(10, 5).Value = "  -0.29%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.160"
$ws.Cells.Item(11, 5).Value = "  +0.05%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "5.04"
$ws.Cells.Item(12, 5).Value = "  -3.35%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.336"
$ws.Cells.Item(13, 5).Value = "  -2.87%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "25.89"
$ws.Cells.Item(14, 5).Value = "  -2.56%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.0000169"
$ws.Cells.Item(15, 5).Value = "  -2.65%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "2.820.17"

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "60.836.24"
$ws.Cells.Item(17, 5).Value = "  -2.06%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "2.418.93"
$ws.Cells.Item(18, 5).Value = "  -0.94%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "10.48"
$ws.Cells.Item(19, 5).Value = "  -4.14%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "7.22"
$ws.Cells.Item(20, 5).Value = "  +0.74%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "320.51"
$ws.Cells.Item(21, 5).Value = "  -3.09%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "4.00"
$ws.Cells.Item(22, 5).Value = "  -2.95%  "

$ws.Cells.Item(24, 5).Value = "  +0.20%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "1.87"
$ws.Cells.Item(25, 5).Value = "  -5.93%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "64.76"
$ws.Cells.Item(26, 5).Value = "  -1.76%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "8.70"
$ws.Cells.Item(27, 5).Value = "  -7.87%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "570.97"
$ws.Cells.Item(28, 5).Value = "  -8.59%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0₃0904"
$ws.Cells.Item(30, 5).Value = "  -5.76%  "

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "7.80"
$ws.Cells.Item(31, 5).Value = "  -2.84%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "1.33"
$ws.Cells.Item(32, 5).Value = "  -8.00%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.82"
$ws.Cells.Item(33, 5).Value = "  -2.86%  "

$ws.Cells.Item(34, 5).Value = "  -6.83%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.01"
$ws.Cells.Item(35, 5).Value = "  +0.58%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "4.58"
$ws.Cells.Item(36, 5).Value = "  -7.37%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.365"
$ws.Cells.Item(37, 5).Value = "  -3.53%  "

$ws.Cells.Item(38, 2).Value = "Monero"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "149.21"
$ws.Cells.Item(38, 5).Value = "  -1.20%  "

$ws.Cells.Item(39, 2).Value = "ImmutableX"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.37"
$ws.Cells.Item(39, 5).Value = "  -5.06%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "18.06"
$ws.Cells.Item(40, 5).Value = "  -1.60%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "5.03"
$ws.Cells.Item(41, 5).Value = "  -5.26%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.00"
$ws.Cells.Item(42, 5).Value = "  +0.06%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "41.72"
$ws.Cells.Item(43, 5).Value = "  -1.76%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.65"
$ws.Cells.Item(44, 5).Value = "  -6.43%  "

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "2.31"
$ws.Cells.Item(45, 5).Value = "  -7.02%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "0.0₆0282"
$ws.Cells.Item(46, 5).Value = "  +17.13%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "139.72"
$ws.Cells.Item(47, 5).Value = "  -2.82%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "3.48"
$ws.Cells.Item(48, 5).Value = "  -4.55%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "0.589"
$ws.Cells.Item(49, 5).Value = "  -2.03%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.0502"
$ws.Cells.Item(50, 5).Value = "  -4.72%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "19.14"
$ws.Cells.Item(51, 5).Value = "  -2.30%  "
